$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the typo'd date "30/02/2025" -> "30/03/2025".
#    Word originally stores this as a single run; the author retyped
#    just the "2" -> "3", which splits the run in three: "30/0" | "3" | "/2025".
# ------------------------------------------------------------------
$dateRange = $d.Content
$dateRange.Find.Execute("30/02/2025") | Out-Null
$dateStart = $dateRange.Start

# Isolate the single digit ("2" at offset 4 of "30/02/2025") into its own
# run *before* touching its text - toggling a character property on/off
# forces a run break without disturbing the rsid of the surrounding runs.
$digitRange = $d.Range($dateStart + 4, $dateStart + 5)
$digitRange.Font.Bold = $true
$digitRange.Font.Bold = $false

# Now retype just that isolated run's content: "2" -> "3".
$digitRange2 = $d.Range($dateStart + 4, $dateStart + 5)
$digitRange2.Find.Execute("2", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "3", 2) | Out-Null

# Retyping re-merged the three pieces back into a single "30/03/2025" run,
# so re-split it into three runs again: "30/0" | "3" | "/2025".
$d.Range($dateStart, $dateStart + 4).Font.Bold = $true
$d.Range($dateStart, $dateStart + 4).Font.Bold = $false
$d.Range($dateStart + 4, $dateStart + 5).Font.Bold = $true
$d.Range($dateStart + 4, $dateStart + 5).Font.Bold = $false
$d.Range($dateStart + 5, $dateStart + 10).Font.Bold = $true
$d.Range($dateStart + 5, $dateStart + 10).Font.Bold = $false

# ------------------------------------------------------------------
# 2) "Ana Karoline, Lays Abreu, Vitor Restini" - the spell checker
#    flags "Lays" and "Restini", which splits the single run into four
#    runs around those two words.
# ------------------------------------------------------------------
$presentesRange = $d.Content
$presentesRange.Find.Execute("Ana Karoline, Lays Abreu, Vitor Restini") | Out-Null
$pStart = $presentesRange.Start

# "Ana Karoline, " | "Lays" | " Abreu, Vitor " | "Restini"
$d.Range($pStart, $pStart + 14).Font.Bold = $true
$d.Range($pStart, $pStart + 14).Font.Bold = $false
$d.Range($pStart + 14, $pStart + 18).Font.Bold = $true
$d.Range($pStart + 14, $pStart + 18).Font.Bold = $false
$d.Range($pStart + 18, $pStart + 32).Font.Bold = $true
$d.Range($pStart + 18, $pStart + 32).Font.Bold = $false
$d.Range($pStart + 32, $pStart + 39).Font.Bold = $true
$d.Range($pStart + 32, $pStart + 39).Font.Bold = $false

# ------------------------------------------------------------------
# 3) "Ana Beatriz Zinatto, Luiz Felipe" - the spell checker flags
#    "Zinatto", which splits the single run into three runs.
# ------------------------------------------------------------------
$ausentesRange = $d.Content
$ausentesRange.Find.Execute("Ana Beatriz Zinatto, Luiz Felipe") | Out-Null
$aStart = $ausentesRange.Start

# "Ana Beatriz " | "Zinatto" | ", Luiz Felipe"
$d.Range($aStart, $aStart + 12).Font.Bold = $true
$d.Range($aStart, $aStart + 12).Font.Bold = $false
$d.Range($aStart + 12, $aStart + 19).Font.Bold = $true
$d.Range($aStart + 12, $aStart + 19).Font.Bold = $false
$d.Range($aStart + 19, $aStart + 32).Font.Bold = $true
$d.Range($aStart + 19, $aStart + 32).Font.Bold = $false
